# Extend "1_auditory_only" Sheet1 from 4 trial-blocks to 16 trial-blocks.
# Rows 2-5 get refreshed timing/ear values (the run was re-timed / the
# trailing "image name" label now reads "none" instead of the old fixed
# image-name token), and rows 6-17 are brand-new blocks (5 through 16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (unchanged labels, rewritten for safety)
$headers = @("run_num", "block_num", "start_time", "play_duration", "ear", "hand")
for ($col = 1; $col -le $headers.Length; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# run_num, block_num, start_time, play_duration, ear, hand
$rows = @(
    @(1, 1,  5.0003231000000001, 17.838872200000001, "L", "none"),
    @(1, 2,  18.173861800000001, 26.3379297,          "R", "none"),
    @(1, 3,  31.000124199999998, 39.2431196,          "R", "none"),
    @(1, 4,  44.000155700000001, 52.021312899999998,  "L", "none"),
    @(1, 5,  57.000289500000001, 65.244836899999996,  "R", "none"),
    @(1, 6,  70.000196799999998, 78.261128299999996,  "R", "none"),
    @(1, 7,  83.000134799999998, 91.178648699999997,  "R", "none"),
    @(1, 8,  96.000094899999993, 104.0448224,         "R", "none"),
    @(1, 9,  109.0002502,        116.9767372,         "L", "none"),
    @(1, 10, 122.000197,         130.0201093,         "L", "none"),
    @(1, 11, 135.0002489,        142.88158000000001,  "R", "none"),
    @(1, 12, 148.00012599999999, 156.1707146,         "L", "none"),
    @(1, 13, 161.00008769999999, 169.1634689,         "L", "none"),
    @(1, 14, 174.00021580000001, 181.8918802,         "R", "none"),
    @(1, 15, 187.00013670000001, 194.96641450000001,  "L", "none"),
    @(1, 16, 200.0001575,        208.17964670000001,  "L", "none")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r++
}

# start_time column now holds longer values -> widen its best-fit column.
$ws.Columns.Item(3).ColumnWidth = 10.75
